$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column header in H1, copying the formatting (bold,
# centered, bordered) from the adjacent existing header cell G1 so the
# new column matches the look of "TB", "d2S", "K", "IP", "Win", "sum".
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data value for the new "Save" column in row 2.
$ws.Range("H2").Value = 0
